$wb = $excel.ActiveWorkbook

# This script applies the data refresh produced by the scheduled Sheets runner
# to the Spriggan_Profits workbook. For each changed leve row, the current market
# price columns (H-N) are updated in place; some rows gain or lose a profit cell
# depending on whether the computed value becomes zero/non-zero.

$ws = $wb.Worksheets.Item("ALC")
# Row 46
$ws.Range("H46").Value = 4500
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()
# Row 60
$ws.Range("H60").Value = 4500
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").ClearContents()
# Row 62
$ws.Range("H62").Value = 5000.5
$ws.Range("J62").Value = 5000.5
$ws.Range("L62").Value = 5000.5
$ws.Range("N62").Value = -6248.5
# Row 64
$ws.Range("H64").Value = 15388807
$ws.Range("I64").Value = 22226166
$ws.Range("K64").Value = 22226166
$ws.Range("M64").Value = -22225918
# Row 65
$ws.Range("H65").Value = 5000.5
$ws.Range("J65").Value = 5000.5
$ws.Range("L65").Value = 25002.5
$ws.Range("N65").Value = -31242.5
# Row 67
$ws.Range("H67").Value = 15388807
$ws.Range("I67").Value = 22226166
$ws.Range("K67").Value = 22226166
$ws.Range("M67").Value = -22225308
# Row 86
$ws.Range("H86").Value = 2192.6667
$ws.Range("I86").Value = 2248.8572
$ws.Range("J86").Value = 1996
$ws.Range("K86").Value = 2248.8572
$ws.Range("L86").Value = 1996
$ws.Range("M86").Value = -1125.8572
$ws.Range("N86").Value = -4242
# Row 89
$ws.Range("H89").Value = 2192.6667
$ws.Range("I89").Value = 2248.8572
$ws.Range("J89").Value = 1996
$ws.Range("K89").Value = 11244.286
$ws.Range("L89").Value = 9980
$ws.Range("M89").Value = -5628.286
$ws.Range("N89").Value = -21212
# Row 106
$ws.Range("H106").Value = 2520.261
$ws.Range("I106").Value = 1970.5
$ws.Range("J106").Value = 4499.4
$ws.Range("K106").Value = 1970.5
$ws.Range("L106").Value = 4499.4
$ws.Range("M106").Value = -1339.5
$ws.Range("N106").Value = -5761.4
# Row 112
$ws.Range("H112").Value = 71817.03
$ws.Range("J112").Value = 58765
$ws.Range("L112").Value = 176295
$ws.Range("N112").Value = -178511

$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 145.45454
$ws.Range("I5").Value = 384.66666
$ws.Range("J5").Value = 55.75
$ws.Range("K5").Value = 384.66666
$ws.Range("L5").Value = 55.75
$ws.Range("M5").Value = -272.66666
$ws.Range("N5").Value = -279.75
# Row 63
$ws.Range("H63").Value = 2695.1667
$ws.Range("I63").Value = 2522.6365
$ws.Range("J63").Value = 4593
$ws.Range("K63").Value = 2522.6365
$ws.Range("L63").Value = 4593
$ws.Range("M63").Value = -1836.6365
$ws.Range("N63").Value = -5965
# Row 66
$ws.Range("H66").Value = 2695.1667
$ws.Range("I66").Value = 2522.6365
$ws.Range("J66").Value = 4593
$ws.Range("K66").Value = 12613.1825
$ws.Range("L66").Value = 22965
$ws.Range("M66").Value = -9181.1825
$ws.Range("N66").Value = -29829
# Row 102
$ws.Range("H102").Value = 10101721
$ws.Range("I102").Value = 11364186
$ws.Range("K102").Value = 11364186
$ws.Range("M102").Value = -11362564

$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 145.45454
$ws.Range("I4").Value = 384.66666
$ws.Range("J4").Value = 55.75
$ws.Range("K4").Value = 384.66666
$ws.Range("L4").Value = 55.75
$ws.Range("M4").Value = -269.66666
$ws.Range("N4").Value = -285.75
# Row 7
$ws.Range("H7").Value = 12965.875
$ws.Range("I7").Value = 14751
$ws.Range("K7").Value = 14751
$ws.Range("M7").Value = -14638
# Row 94
$ws.Range("H94").Value = 1224.75
$ws.Range("I94").Value = 1300
$ws.Range("J94").Value = 999
$ws.Range("K94").Value = 1300
$ws.Range("L94").Value = 999
$ws.Range("M94").Value = -849
$ws.Range("N94").Value = -1901
# Row 105
$ws.Range("H105").Value = 2882.889
$ws.Range("I105").Value = 2123.7
$ws.Range("J105").Value = 3831.875
$ws.Range("K105").Value = 2123.7
$ws.Range("L105").Value = 3831.875
$ws.Range("M105").Value = -376.6999999999998
$ws.Range("N105").Value = -7325.875

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 17074.666
$ws.Range("I22").Value = 33549.668
$ws.Range("J22").Value = 599.6667
$ws.Range("K22").Value = 33549.668
$ws.Range("L22").Value = 599.6667
$ws.Range("M22").Value = -33199.668
$ws.Range("N22").Value = -1299.6667
# Row 68
$ws.Range("H68").Value = 100000
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
# Row 71
$ws.Range("H71").Value = 100000
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
# Row 86
$ws.Range("H86").Value = 4770.3335
$ws.Range("J86").Value = 4441.5
$ws.Range("L86").Value = 4441.5
$ws.Range("N86").Value = -6687.5
# Row 89
$ws.Range("H89").Value = 4770.3335
$ws.Range("J89").Value = 4441.5
$ws.Range("L89").Value = 22207.5
$ws.Range("N89").Value = -33439.5
# Row 122
$ws.Range("H122").Value = 3385.2
$ws.Range("J122").Value = 5300
$ws.Range("L122").Value = 15900
$ws.Range("N122").Value = -20800
# Row 134
$ws.Range("H134").Value = 31253030
$ws.Range("I134").Value = 35717256
$ws.Range("K134").Value = 107151768
$ws.Range("M134").Value = -107149233

$ws = $wb.Worksheets.Item("CUL")
# Row 124
$ws.Range("H124").Value = 500
$ws.Range("I124").Value = 500
$ws.Range("K124").Value = 1500
$ws.Range("M124").Value = 3410
# Row 126
$ws.Range("H126").Value = 9896
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 9896
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 29688
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -39568
# Row 132
$ws.Range("H132").Value = 1702.6923
$ws.Range("I132").Value = 1695.1666
$ws.Range("K132").Value = 15256.4994
$ws.Range("M132").Value = -12726.4994
# Row 141
$ws.Range("H141").Value = 2007.5
$ws.Range("I141").Value = 2007.5
$ws.Range("K141").Value = 6022.5
$ws.Range("M141").Value = -842.5

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 772.8333
$ws.Range("I46").Value = 772.8333
$ws.Range("K46").Value = 772.8333
$ws.Range("M46").Value = -584.8333
# Row 82
$ws.Range("H82").Value = 1699.5
$ws.Range("I82").Value = 2400
$ws.Range("K82").Value = 2400
$ws.Range("M82").Value = -2039
# Row 85
$ws.Range("H85").Value = 1699.5
$ws.Range("I85").Value = 2400
$ws.Range("K85").Value = 2400
$ws.Range("M85").Value = -1152
# Row 111
$ws.Range("H111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 54
$ws.Range("H54").Value = 10927.667
$ws.Range("I54").Value = 8891.5
$ws.Range("J54").Value = 15000
$ws.Range("K54").Value = 8891.5
$ws.Range("L54").Value = 15000
$ws.Range("M54").Value = -8371.5
$ws.Range("N54").Value = -16040
# Row 62
$ws.Range("H62").Value = 23000
$ws.Range("J62").Value = 23000
$ws.Range("L62").Value = 23000
$ws.Range("N62").Value = -24248
# Row 65
$ws.Range("H65").Value = 23000
$ws.Range("J65").Value = 23000
$ws.Range("L65").Value = 115000
$ws.Range("N65").Value = -121240
# Row 140
$ws.Range("H140").Value = 72412.836
$ws.Range("J140").Value = 72412.836
$ws.Range("L140").Value = 72412.836
$ws.Range("N140").Value = -82772.836
# Row 141
$ws.Range("H141").Value = 75995.664
$ws.Range("J141").Value = 75794.8
$ws.Range("L141").Value = 75794.8
$ws.Range("N141").Value = -86154.8

Write-Host "Spriggan_Profits leve-price refresh applied"